# inspector update and camera features
#
# - B2 / C2 dialogue script text updated (new SAY beats + new dialogue line
#   with renamed commands and extra panel lines)
# - selection cursor moved from C2 to C8
# - sheet default column width nudged (11.5703125 -> 11.58984375)
# - row 2 height grown from 53.7 to 99.3 (taller text needs more room)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dialogueCell = @'
[ ["SAY" ,{"name": "dahlia", "mood": "happy"}], ["SAY" ,{"name": "dahlia", "mood": "happy"}], ["SAY" ,{"name": "dahlia", "mood": "thinking"}], ["SAY" ,{"name": "dahlia", "mood": "happy"}] ,["SAY",{"name":"???"}]]
'@

$textCell = @'
["#set_speed;5##set_track;dahlia#Hello everyone :D this is a long text that will likely overflow of this message box, spilling it’s contents outside and litter the universe.","#set_speed;4#However as you can see, this panel is adapting to fit everything :D even with long-ass words, watch : thisisaverylongwordandidon’tknowwhattowritetomakeitlongerhopefullyyougetitbynowherehavesomegibberishforgoodmeasuregbnhirltuhgtbiuthvbiliu","hmmm…","Nice weather we’re having today, isn’t it ?","#set_speed;0.2#yeah it’s nice.#set_delay;1000#"]
'@

$ws.Range("B2").Value = $dialogueCell
$ws.Range("C2").Value = $textCell

# Row 2 grew taller to fit the longer dialogue text.
$ws.Rows.Item(2).RowHeight = 99.3

# Default sheet column width tweak.
$ws.StandardWidth = 11.58984375

# Move the live selection cursor to C8.
$null = $ws.Range("C8").Select()
